$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 40: fill in action (D40) and date (E40) ---
$ws.Range("D40").Value = "done"
$ws.Range("E40").Value2 = 45985

# --- Add new rows 41 and 42, copying formatting (styles) from row 40 ---
$ws.Range("A40:E40").Copy()
$ws.Range("A41:E42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 41
$ws.Range("A41").Value = "make meglm work without rmcoll"
$ws.Range("B41").Value = "Ian"
$ws.Range("C41").Value2 = 45985
$ws.Range("D41").Value = ""
$ws.Range("E41").Value = ""

# Row 42
$ws.Range("A42").Value = "syntax 2: expand 1st derivative to act after any collinearity"
$ws.Range("B42").Value = "Ian"
$ws.Range("C42").Value2 = 45985
$ws.Range("D42").Value = ""
$ws.Range("E42").Value = ""

# --- Resize the Excel table (ListObject) to include the two new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:E42"))

# --- Update selection/view to match end state shown in diff ---
$ws.Range("B42:C42").Select()
$excel.ActiveWindow.ScrollRow = 33
